$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 9666
$ws.Range("J44").Value = 9666
$ws.Range("L44").Value = 9666
$ws.Range("N44").Value = -10590
$ws.Range("H80").Value = 4167.9165
$ws.Range("I80").Value = 3655
$ws.Range("J80").Value = 4680.8335
$ws.Range("K80").Value = 10965
$ws.Range("L80").Value = 14042.5005
$ws.Range("M80").Value = -9967
$ws.Range("N80").Value = -16038.5005
$ws.Range("H83").Value = 4167.9165
$ws.Range("I83").Value = 3655
$ws.Range("J83").Value = 4680.8335
$ws.Range("K83").Value = 32895
$ws.Range("L83").Value = 42127.5015
$ws.Range("M83").Value = -27903
$ws.Range("N83").Value = -52111.5015
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H98").Value = 1232.7142
$ws.Range("I98").Value = 1104.8334
$ws.Range("K98").Value = 1104.8334
$ws.Range("M98").Value = 393.1666
$ws.Range("H101").Value = 25001500
$ws.Range("I101").Value = 33333666
$ws.Range("K101").Value = 100000998
$ws.Range("M101").Value = -99999376
$ws.Range("H104").Value = 131.66667
$ws.Range("I104").Value = 131.66667
$ws.Range("K104").Value = 395.00001
$ws.Range("M104").Value = 1351.99999
$ws.Range("H116").Value = 1999.3334
$ws.Range("J116").Value = 1999.3334
$ws.Range("L116").Value = 1999.3334
$ws.Range("N116").Value = -8883.3334
$ws.Range("H122").Value = 1232.7142
$ws.Range("I122").Value = 1104.8334
$ws.Range("K122").Value = 3314.5002
$ws.Range("M122").Value = -864.5001999999999
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2574.6667
$ws.Range("I61").Value = 2769.9
$ws.Range("J61").Value = 1598.5
$ws.Range("K61").Value = 2769.9
$ws.Range("L61").Value = 1598.5
$ws.Range("M61").Value = -2557.9
$ws.Range("N61").Value = -2022.5
$ws.Range("H74").Value = 638.35297
$ws.Range("I74").Value = 522
$ws.Range("J74").Value = 2500
$ws.Range("K74").Value = 522
$ws.Range("L74").Value = 2500
$ws.Range("M74").Value = 352
$ws.Range("N74").Value = -4248
$ws.Range("H77").Value = 638.35297
$ws.Range("I77").Value = 522
$ws.Range("J77").Value = 2500
$ws.Range("K77").Value = 2610
$ws.Range("L77").Value = 12500
$ws.Range("M77").Value = 1758
$ws.Range("N77").Value = -21236
$ws.Range("H88").Value = 3999.6667
$ws.Range("J88").Value = 3999.75
$ws.Range("L88").Value = 3999.75
$ws.Range("N88").Value = -4811.75
$ws.Range("H91").Value = 3999.6667
$ws.Range("J91").Value = 3999.75
$ws.Range("L91").Value = 3999.75
$ws.Range("N91").Value = -6807.75
$ws.Range("H97").Value = 651.6667
$ws.Range("I97").Value = 629.3333
$ws.Range("K97").Value = 629.3333
$ws.Range("M97").Value = -133.3333
$ws.Range("H136").Value = 2574.6667
$ws.Range("I136").Value = 2769.9
$ws.Range("J136").Value = 1598.5
$ws.Range("K136").Value = 8309.700000000001
$ws.Range("L136").Value = 4795.5
$ws.Range("M136").Value = -5759.700000000001
$ws.Range("N136").Value = -9895.5
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1143.4
$ws.Range("I99").Value = 929.25
$ws.Range("K99").Value = 929.25
$ws.Range("M99").Value = 568.75
$ws.Range("H134").Value = 2946.3333
$ws.Range("I134").Value = 2946.3333
$ws.Range("K134").Value = 8838.999899999999
$ws.Range("M134").Value = -6303.999899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1349.5555
$ws.Range("I31").Value = 1092.2858
$ws.Range("K31").Value = 1092.2858
$ws.Range("M31").Value = -797.2858000000001
$ws.Range("H34").Value = 1349.5555
$ws.Range("I34").Value = 1092.2858
$ws.Range("K34").Value = 1092.2858
$ws.Range("M34").Value = -890.2858000000001
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H58").Value = 2414.4546
$ws.Range("I58").Value = 2172.6
$ws.Range("K58").Value = 2172.6
$ws.Range("M58").Value = -1969.6
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H99").Value = 3750
$ws.Range("I99").Value = 2500
$ws.Range("K99").Value = 2500
$ws.Range("M99").Value = -1002
$ws.Range("H126").Value = 3750
$ws.Range("I126").Value = 2500
$ws.Range("K126").Value = 7500
$ws.Range("M126").Value = -5030
$ws.Range("H136").Value = 2414.4546
$ws.Range("I136").Value = 2172.6
$ws.Range("K136").Value = 6517.799999999999
$ws.Range("M136").Value = -3967.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 1000
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H58").Value = 2025
$ws.Range("I58").Value = 2025
$ws.Range("K58").Value = 6075
$ws.Range("M58").Value = -5947
$ws.Range("H114").Value = 2289.3125
$ws.Range("I114").Value = 1506.8572
$ws.Range("J114").Value = 2897.889
$ws.Range("K114").Value = 4520.571599999999
$ws.Range("L114").Value = 8693.667000000001
$ws.Range("M114").Value = -1266.571599999999
$ws.Range("N114").Value = -15201.667
$ws.Range("H141").Value = 7205.9
$ws.Range("I141").Value = 7117.6665
$ws.Range("K141").Value = 21352.9995
$ws.Range("M141").Value = -16172.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 52371.668
$ws.Range("J15").Value = 52371.668
$ws.Range("L15").Value = 52371.668
$ws.Range("N15").Value = -52947.668
$ws.Range("H80").Value = 4416.6665
$ws.Range("I80").Value = 5000
$ws.Range("J80").Value = 4300
$ws.Range("K80").Value = 5000
$ws.Range("L80").Value = 4300
$ws.Range("M80").Value = -4002
$ws.Range("N80").Value = -6296
$ws.Range("H81").Value = 52371.668
$ws.Range("J81").Value = 52371.668
$ws.Range("L81").Value = 52371.668
$ws.Range("N81").Value = -54367.668
$ws.Range("H83").Value = 4416.6665
$ws.Range("I83").Value = 5000
$ws.Range("J83").Value = 4300
$ws.Range("K83").Value = 25000
$ws.Range("L83").Value = 21500
$ws.Range("M83").Value = -20008
$ws.Range("N83").Value = -31484
$ws.Range("H84").Value = 52371.668
$ws.Range("J84").Value = 52371.668
$ws.Range("L84").Value = 157115.004
$ws.Range("N84").Value = -167099.004
$ws.Range("H139").Value = 50000
$ws.Range("J139").Value = 50000
$ws.Range("L139").Value = 50000
$ws.Range("N139").Value = -60280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 9840.143
$ws.Range("I132").Value = 14343.75
$ws.Range("K132").Value = 43031.25
$ws.Range("M132").Value = -40501.25
$ws.Range("H136").Value = 1766.0834
$ws.Range("I136").Value = 1290.2727
$ws.Range("K136").Value = 3870.8181
$ws.Range("M136").Value = -1320.8181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3333859
$ws.Range("I81").Value = 787.5
$ws.Range("K81").Value = 1575
$ws.Range("M81").Value = -514
$ws.Range("H84").Value = 3333859
$ws.Range("I84").Value = 787.5
$ws.Range("K84").Value = 7875
$ws.Range("M84").Value = -2571
$ws.Range("H94").Value = 375162.5
$ws.Range("J94").Value = 375162.5
$ws.Range("L94").Value = 375162.5
$ws.Range("N94").Value = -376964.5
$ws.Range("H132").Value = 2328.4285
$ws.Range("J132").Value = 1175.5
$ws.Range("L132").Value = 3526.5
$ws.Range("N132").Value = -8586.5
